$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.199.38"
$ws.Range("E2").Value = "  +3.45%  "
$ws.Range("D3").Value = "2.623.23"
$ws.Range("E3").Value = "  +1.75%  "
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.ClearFormats()
$ws.Range("E4").Value = "  +0.25%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "568.09"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +6.05%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "145.32"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +2.49%  "
$ws.Range("E7").Value = "  -0.15%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.609"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +4.43%  "
$ws.Range("D9").Value = "2.642.19"
$ws.Range("E9").Value = "  +2.43%  "
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("E11").Value = "  +5.30%  "
$ws.Range("E12").Value = "  +7.48%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.344"
$cell.ClearFormats()
$ws.Range("E13").Value = "  +3.91%  "
$ws.Range("D14").Value = "3.091.61"
$ws.Range("E14").Value = "  +1.94%  "
$ws.Range("D15").Value = "60.173.59"
$ws.Range("E15").Value = "  +3.53%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "22.01"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +6.70%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.0000138"
$cell.ClearFormats()
$ws.Range("E17").Value = "  +4.23%  "
$ws.Range("D18").Value = "2.639.37"
$ws.Range("E18").Value = "  +2.72%  "
$ws.Range("E19").Value = "  +2.54%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "341.13"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +1.96%  "
$ws.Range("E21").Value = "  +3.65%  "
$ws.Range("E22").Value = "  +3.75%  "
$ws.Range("E23").Value = "  -0.11%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "65.95"
$cell.ClearFormats()
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("E25").Value = "  +5.17%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.163"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +3.28%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +0.20%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "7.38"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +4.82%  "
$ws.Range("D29").Value = "0.0₃0801"
$ws.Range("E29").Value = "  +10.65%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  +4.56%  "
$ws.Range("E32").Value = "  +5.00%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "158.32"
$cell.ClearFormats()
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("E34").Value = "  +1.43%  "
$ws.Range("E35").Value = "  +5.64%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.15"
$cell.ClearFormats()
$ws.Range("E36").Value = "  +5.24%  "
$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.882"
$cell.ClearFormats()
$ws.Range("E37").Value = "  +6.55%  "
$ws.Range("E38").Value = "  +8.27%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "37.50"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("E40").Value = "  +7.15%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "298.80"
$cell.ClearFormats()
$ws.Range("E41").Value = "  +5.33%  "
$ws.Range("E42").Value = "  +1.47%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.995"
$cell.ClearFormats()
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("E44").Value = "  +4.05%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.602"
$cell.ClearFormats()
$ws.Range("E45").Value = "  +2.37%  "
$ws.Range("E46").Value = "  +2.13%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "128.03"
$cell.ClearFormats()
$ws.Range("E47").Value = "  +15.46%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "19.36"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +5.71%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "10.69"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("E50").Value = "  +3.83%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "4.65"
$cell.ClearFormats()
$ws.Range("E51").Value = "  +6.54%  "
